$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that already contain a MultivariableGLS (column H) value and must stay untouched
$skipRows = @(14, 28)

for ($r = 2; $r -le 27; $r++) {
    if ($skipRows -contains $r) {
        continue
    }
    $ws.Cells.Item($r, 8).Value = 1
}
